$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 15 ("Giai thuat SPIMI") - Content Placeholder 2 (Shapes.Item(3))
# ---------------------------------------------------------------------------
$s15 = $p.Slides.Item(15)
$content = $s15.Shapes.Item(3)
$ctr = $content.TextFrame.TextRange

# Paragraph 2: "Khong can cung cap ma tu duy nhat tren toan bo du lieu;"
#           -> "Khong can quan ly ma tu trong qua trinh xay dung chi muc khoi;"
$ctr.Paragraphs(2).Text = "Không cần quản lý mã từ trong quá trình xây dựng chỉ mục khối;"

# Paragraph 3: "Khong can luu tu dien day du cho bo du lieu trong bo nho."
#           -> "Kich thuoc khoi trong bo nho lon hon so voi BSBI."
$ctr.Paragraphs(3).Text = "Kích thước khối trong bộ nhớ lớn hơn so với BSBI."

# Paragraph 5: drop the trailing period
$ctr.Paragraphs(5).Text = "Không cần thực hiện sắp xếp danh sách thẻ định vị"

# Add a brand-new paragraph 6 at the end of this placeholder.
$content.TextFrame.TextRange.InsertAfter("`rTiết kiệm bộ nhớ hơn so với BSBI.") | Out-Null

Write-Host "slide15 content placeholder done"

# ---------------------------------------------------------------------------
# Slide 15 - "TextBox 1" (Shapes.Item(5)): rewrite both bullet lines, drop
# the 3rd paragraph (it gets folded into paragraph 2's tail), bump the font
# size 20 -> 22pt and grow the box to fit.
# ---------------------------------------------------------------------------
$tb1 = $s15.Shapes.Item(5)
$tbtr = $tb1.TextFrame.TextRange
$tbtr.Text = "Xây dựng chỉ mục một lượt trong bộ nhớ : SPIMI: Single-pass in-memory indexing;`rXây dựng chỉ mục ngược đầy đủ cho mỗi khối -> Sắp xếp từ điển cục bộ -> Ghi ra đĩa -> hợp nhất khối"
$tbtr.Font.Size = 22
$tb1.Height = 113.90161

Write-Host "slide15 textbox1 done"
